$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column G ("Line / Loại lý do") ---
# This shifts the old G:K (Đơn vị-data .. Trung tâm chi phí) one column
# to the right, becoming H:L, and leaves the new G column blank.
$ws.Columns("G:G").Insert()

# --- Helper to write a header cell whose trailing "(*)" marker is bold/red ---
function Set-HeaderStar {
    param($cellRef, $text, $starStart, $starLen)
    $ws.Range($cellRef).Value = $text
    $chars = $ws.Range($cellRef).Characters($starStart, $starLen)
    $chars.Font.Bold = $true
    $chars.Font.Color = 255
    $chars.Font.Name = "Calibri"
    $chars.Font.Size = 11
}

Set-HeaderStar "A1" "Ngày kế hoạch (*)"        15 3
Set-HeaderStar "B1" "Loại phiếu (*)"           12 3
Set-HeaderStar "C1" "Line / Sản phẩm (*)"      17 3
Set-HeaderStar "D1" "Line / Ngày dự kiến (*)"  21 3
Set-HeaderStar "E1" "Line / Số lượng (*)"      17 3
Set-HeaderStar "F1" "Line / Đơn vị (*)"        15 3
Set-HeaderStar "G1" "Line / Loại lý do (*)"    18 4
Set-HeaderStar "H1" "Line / Lý do nhập (*)"    19 3
Set-HeaderStar "I1" "Line / Đến kho (*)"       15 4

# --- Re-fit the columns whose header text grew because of the new "(*)" marker ---
$ws.Columns("D:D").ColumnWidth = 19.666666666666668
$ws.Columns("E:E").ColumnWidth = 16
$ws.Columns("F:F").ColumnWidth = 13.833333333333334
$ws.Columns("G:G").ColumnWidth = 16.333333333333332
$ws.Columns("I:I").ColumnWidth = 15.333333333333334

# --- Move the active selection, matching the saved cursor position ---
$ws.Range("G10").Select()

# --- Page setup: force portrait orientation ---
$ws.PageSetup.Orientation = 1
